$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.055.55'
$ws.Range('E2').Value = '  +1.83%  '
$ws.Range('D3').Value = '3.931.34'
$ws.Range('E3').Value = '  +2.41%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = '''609.93'
$ws.Range('E5').Value = '  +1.48%  '
$ws.Range('D6').Value = '''170.62'
$ws.Range('E6').Value = '  +5.50%  '
$ws.Range('D7').Value = '3.932.02'
$ws.Range('E7').Value = '  +2.49%  '
$ws.Range('E8').Value = '  +0.16%  '
$ws.Range('E9').Value = '  +1.26%  '
$ws.Range('E10').Value = '  +1.30%  '
$ws.Range('D11').Value = '''6.43'
$ws.Range('E11').Value = '  +2.10%  '
$ws.Range('D12').Value = '''0.470'
$ws.Range('E12').Value = '  +2.57%  '
$ws.Range('B13').Value = 'ShibaInu'
$ws.Range('C13').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D13').Value = '''0.0000257'
$ws.Range('E13').Value = '  +5.54%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').Value = '''38.51'
$ws.Range('E14').Value = '  +4.58%  '
$ws.Range('D15').Value = '4.596.70'
$ws.Range('E15').Value = '  +2.48%  '
$ws.Range('D16').Value = '3.906.74'
$ws.Range('E16').Value = '  +2.52%  '
$ws.Range('D17').Value = '70.098.76'
$ws.Range('E17').Value = '  +1.73%  '
$ws.Range('D18').Value = '''18.66'
$ws.Range('E18').Value = '  +9.08%  '
$ws.Range('D19').Value = '''7.65'
$ws.Range('E19').Value = '  +1.87%  '
$ws.Range('D21').Value = '''11.17'
$ws.Range('E21').Value = '  -1.20%  '
$ws.Range('D22').Value = '''494.64'
$ws.Range('E22').Value = '  +2.17%  '
$ws.Range('D23').Value = '''0.749'
$ws.Range('E23').Value = '  +4.38%  '
$ws.Range('E24').Value = '  +5.40%  '
$ws.Range('D25').Value = '''86.15'
$ws.Range('E25').Value = '  +2.56%  '
$ws.Range('E26').Value = '  +2.40%  '
$ws.Range('D27').Value = '''12.33'
$ws.Range('E27').Value = '  +2.09%  '
$ws.Range('E28').Value = '  +1.93%  '
$ws.Range('E29').Value = '  +0.16%  '
$ws.Range('D30').Value = '''3.00'
$ws.Range('E30').Value = '  +1.72%  '
$ws.Range('E31').Value = '  +3.53%  '
$ws.Range('D32').Value = '4.084.10'
$ws.Range('E32').Value = '  +2.39%  '
$ws.Range('D33').Value = '''7.87'
$ws.Range('E33').Value = '  -0.55%  '
$ws.Range('D34').Value = '''32.29'
$ws.Range('E34').Value = '  +0.70%  '
$ws.Range('D35').Value = '3.896.61'
$ws.Range('E35').Value = '  +2.91%  '
$ws.Range('E36').Value = '  +1.65%  '
$ws.Range('D37').Value = '''6.16'
$ws.Range('E37').Value = '  +4.58%  '
$ws.Range('E38').Value = '  +2.11%  '
$ws.Range('D39').Value = '''0.142'
$ws.Range('E39').Value = '  +1.67%  '
$ws.Range('E40').Value = '  +12.10%  '
$ws.Range('D41').Value = '''0.330'
$ws.Range('E41').Value = '  +3.61%  '
$ws.Range('E42').Value = '  +0.12%  '
$ws.Range('E43').Value = '  +8.28%  '
$ws.Range('D44').Value = '''439.71'
$ws.Range('E44').Value = '  +0.52%  '
$ws.Range('D45').Value = '''48.39'
$ws.Range('E45').Value = '  -0.21%  '
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').Value = '''0.000279'
$ws.Range('E48').Value = '  +24.61%  '
$ws.Range('E49').Value = '  +3.44%  '
$ws.Range('D50').Value = '''40.84'
$ws.Range('E50').Value = '  +6.05%  '
$ws.Range('D51').Value = '''143.37'
$ws.Range('E51').Value = '  -0.07%  '
